$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pl_mw values for the 380 kV case (rows 2-25, columns B,C,E,F,G,I,K,M,N)
# Row 2
$ws.Range("B2").Value = 0.4756736078547874
$ws.Range("C2").Value = 0.1289437545541432
$ws.Range("E2").Value = 0.1177529936711394
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002504701607980413
$ws.Range("I2").Value = 1.093206866784406
$ws.Range("K2").Value = 0.593591767820925
$ws.Range("M2").Value = 0.2976026010178501
$ws.Range("N2").Value = 2.309051508517086

# Row 3
$ws.Range("B3").Value = 0.4407338459475056
$ws.Range("C3").Value = 0.1197948262781239
$ws.Range("E3").Value = 0.1077341651412667
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002508463030738485
$ws.Range("I3").Value = 1.09187812045267
$ws.Range("K3").Value = 0.5503887842177733
$ws.Range("M3").Value = 0.2743082962894903
$ws.Range("N3").Value = 2.323885184752157

# Row 4
$ws.Range("B4").Value = 0.4195348428377486
$ws.Range("C4").Value = 0.1142564059016564
$ws.Range("E4").Value = 0.1016512026924019
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002510893081447891
$ws.Range("I4").Value = 1.091571430499634
$ws.Range("K4").Value = 0.5241913622026857
$ws.Range("M4").Value = 0.2601746139530121
$ws.Range("N4").Value = 2.333688740256903

# Row 5
$ws.Range("B5").Value = 0.4109598759414155
$ws.Range("C5").Value = 0.1120191538246615
$ws.Range("E5").Value = 0.09918936351323993
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002511913751961539
$ws.Range("I5").Value = 1.091574231279026
$ws.Range("K5").Value = 0.5135982109030977
$ws.Range("M5").Value = 0.2544571973505896
$ws.Range("N5").Value = 2.337858388582681

# Row 6
$ws.Range("B6").Value = 0.4095398594659798
$ws.Range("C6").Value = 0.1116488443588537
$ws.Range("E6").Value = 0.09878159775303885
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002512085072821413
$ws.Range("I6").Value = 1.091582406528573
$ws.Range("K6").Value = 0.5118442006046848
$ws.Range("M6").Value = 0.253510361675211
$ws.Range("N6").Value = 2.338561293436904

# Row 7
$ws.Range("B7").Value = 0.4194189396575325
$ws.Range("C7").Value = 0.1142261540137071
$ws.Range("E7").Value = 0.1016179328831797
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.00251090672338831
$ws.Range("I7").Value = 1.091570951240641
$ws.Range("K7").Value = 0.5240481655894484
$ws.Range("M7").Value = 0.2600973365744608
$ws.Range("N7").Value = 2.333744266973341

# Row 8
$ws.Range("B8").Value = 0.4635735496603957
$ws.Range("C8").Value = 0.1257726833198376
$ws.Range("E8").Value = 0.1142841062603495
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002505973591694908
$ws.Range("I8").Value = 1.092642881019927
$ws.Range("K8").Value = 0.5786268093258968
$ws.Range("M8").Value = 0.2895354049561973
$ws.Range("N8").Value = 2.314021604453615

# Row 9
$ws.Range("B9").Value = 0.5521872955137042
$ws.Range("C9").Value = 0.1490521323495386
$ws.Range("E9").Value = 0.1396793136090793
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002497251509569889
$ws.Range("I9").Value = 1.098797808077002
$ws.Range("K9").Value = 0.6882894195481128
$ws.Range("M9").Value = 0.3486243281278973
$ws.Range("N9").Value = 2.280878319369677

# Row 10
$ws.Range("B10").Value = 0.6185487838850463
$ws.Range("C10").Value = 0.1665586593999819
$ws.Range("E10").Value = 0.1586955540656092
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002491417283121351
$ws.Range("I10").Value = 1.10581042140447
$ws.Range("K10").Value = 0.7705018528452854
$ws.Range("M10").Value = 0.3928977782855654
$ws.Range("N10").Value = 2.259919254882078

# Row 11
$ws.Range("B11").Value = 0.6490164193786256
$ws.Range("C11").Value = 0.1746137200145483
$ws.Range("E11").Value = 0.1674285617553721
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002488886404433843
$ws.Range("I11").Value = 1.109545909393375
$ws.Range("K11").Value = 0.8082679656387768
$ws.Range("M11").Value = 0.4132328736127491
$ws.Range("N11").Value = 2.251124875137734

# Row 12
$ws.Range("B12").Value = 0.6605941699863251
$ws.Range("C12").Value = 0.1776773190529468
$ws.Range("E12").Value = 0.1707476796091143
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.00248794563107518
$ws.Range("I12").Value = 1.111039204650851
$ws.Range("K12").Value = 0.8226223429366826
$ws.Range("M12").Value = 0.4209617436120681
$ws.Range("N12").Value = 2.247901442657721

# Row 13
$ws.Range("B13").Value = 0.6580988989160517
$ws.Range("C13").Value = 0.1770169224561755
$ws.Range("E13").Value = 0.1700323050272772
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002488147461453094
$ws.Range("I13").Value = 1.110714089412184
$ws.Range("K13").Value = 0.8195285004289872
$ws.Range("M13").Value = 0.4192959233845244
$ws.Range("N13").Value = 2.248590909810119

# Row 14
$ws.Range("B14").Value = 0.6499681190759645
$ws.Range("C14").Value = 0.17486549599505
$ws.Range("E14").Value = 0.1677013832030099
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002488808653779438
$ws.Range("I14").Value = 1.109667183656057
$ws.Range("K14").Value = 0.8094478412172634
$ws.Range("M14").Value = 0.4138681607723811
$ws.Range("N14").Value = 2.250857538519881

# Row 15
$ws.Range("B15").Value = 0.6449930362886107
$ws.Range("C15").Value = 0.1735494261178587
$ws.Range("E15").Value = 0.166275212222466
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002489215944702432
$ws.Range("I15").Value = 1.109036188312423
$ws.Range("K15").Value = 0.8032800814227414
$ws.Range("M15").Value = 0.4105472109300479
$ws.Range("N15").Value = 2.252259835506507

# Row 16
$ws.Range("B16").Value = 0.6165632869052615
$ws.Range("C16").Value = 0.1660340955691879
$ws.Range("E16").Value = 0.1581265126510303
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002491585151931216
$ws.Range("I16").Value = 1.105577300634678
$ws.Range("K16").Value = 0.7680411692857092
$ws.Range("M16").Value = 0.3915727864723095
$ws.Range("N16").Value = 2.260508912823298

# Row 17
$ws.Range("B17").Value = 0.5991942256592893
$ws.Range("C17").Value = 0.1614471977312917
$ws.Range("E17").Value = 0.153148847973462
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002493070057788715
$ws.Range("I17").Value = 1.103595303305077
$ws.Range("K17").Value = 0.7465175557241821
$ws.Range("M17").Value = 0.3799827793304473
$ws.Range("N17").Value = 2.265759267738488

# Row 18
$ws.Range("B18").Value = 0.5892302970626133
$ws.Range("C18").Value = 0.1588175125375813
$ws.Range("E18").Value = 0.1502935691236189
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002493935732289764
$ws.Range("I18").Value = 1.102506628689611
$ws.Range("K18").Value = 0.7341722893380052
$ws.Range("M18").Value = 0.3733348334663162
$ws.Range("N18").Value = 2.268848761559397

# Row 19
$ws.Range("B19").Value = 0.5858611944739778
$ws.Range("C19").Value = 0.157928613377976
$ws.Range("E19").Value = 0.1493281424657553
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002494230829080226
$ws.Range("I19").Value = 1.102146826469266
$ws.Range("K19").Value = 0.7299983169272082
$ws.Range("M19").Value = 0.3710870879239749
$ws.Range("N19").Value = 2.269906757877138

# Row 20
$ws.Range("B20").Value = 0.6010404693995781
$ws.Range("C20").Value = 0.1619345916940631
$ws.Range("E20").Value = 0.1536779256773713
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002492910787691206
$ws.Range("I20").Value = 1.103800976810525
$ws.Range("K20").Value = 0.7488052030937524
$ws.Range("M20").Value = 0.3812146565698171
$ws.Range("N20").Value = 2.265193149330628

# Row 21
$ws.Range("B21").Value = 0.6523552309035097
$ws.Range("C21").Value = 0.1754970589250604
$ws.Range("E21").Value = 0.1683857004606892
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002488613968340569
$ws.Range("I21").Value = 1.109972545536245
$ws.Range("K21").Value = 0.8124073292809726
$ws.Range("M21").Value = 0.4154616520563437
$ws.Range("N21").Value = 2.250188872788286

# Row 22
$ws.Range("B22").Value = 0.6861274646909692
$ws.Range("C22").Value = 0.1844386629292671
$ws.Range("E22").Value = 0.1780688354594488
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002485908383672603
$ws.Range("I22").Value = 1.114465170288376
$ws.Range("K22").Value = 0.8542850820933268
$ws.Range("M22").Value = 0.4380098256754081
$ws.Range("N22").Value = 2.24100545911709

# Row 23
$ws.Range("B23").Value = 0.6680810339712195
$ws.Range("C23").Value = 0.1796591799408702
$ws.Range("E23").Value = 0.1728942019470381
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002487343044139811
$ws.Range("I23").Value = 1.112025254109398
$ws.Range("K23").Value = 0.8319056442613828
$ws.Range("M23").Value = 0.4259601428883428
$ws.Range("N23").Value = 2.245849705808467

# Row 24
$ws.Range("B24").Value = 0.6002057147857158
$ws.Range("C24").Value = 0.1617142179009647
$ws.Range("E24").Value = 0.1534387095914198
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002492982756265623
$ws.Range("I24").Value = 1.103707833584572
$ws.Range("K24").Value = 0.7477708675856718
$ws.Range("M24").Value = 0.3806576772367549
$ws.Range("N24").Value = 2.265448870318494

# Row 25
$ws.Range("B25").Value = 0.5279957600931198
$ws.Range("C25").Value = 0.1426847203324542
$ws.Range("E25").Value = 0.1327476745306768
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002499509827420079
$ws.Range("I25").Value = 1.096696811381868
$ws.Range("K25").Value = 0.6583369306146096
$ws.Range("M25").Value = 0.3324903051167567
$ws.Range("N25").Value = 2.289250254603985

